$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "43.044.18"
$ws.Range("E2").Value = "  +0.63%  "
$ws.Range("D3").Value = "2.305.55"
$ws.Range("E3").Value = "  +0.70%  "
$ws.Range("E4").Value = "  -0.19%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "305.47"
$ws.Range("E5").Value = "  +1.85%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "97.42"
$ws.Range("E6").Value = "  +0.97%  "
$ws.Range("E7").Value = "  -1.64%  "
$ws.Range("E8").Value = "  -0.14%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.501"
$ws.Range("E9").Value = "  -0.52%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "35.52"
$ws.Range("E10").Value = "  +0.04%  "
$ws.Range("E11").Value = "  +0.45%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "18.67"
$ws.Range("E12").Value = "  +5.52%  "
$ws.Range("E13").Value = "  +1.71%  "
$ws.Range("E14").Value = "  +2.46%  "
$ws.Range("D15").Value = "2.664.48"
$ws.Range("E15").Value = "  +0.66%  "
$ws.Range("D16").Value = "2.309.63"
$ws.Range("E16").Value = "  +1.15%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.784"
$ws.Range("E17").Value = "  +1.49%  "
$ws.Range("D18").Value = "42.918.53"
$ws.Range("E18").Value = "  +0.50%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "12.61"
$ws.Range("E19").Value = "  -1.09%  "
$ws.Range("E20").Value = "  -0.31%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.06"
$ws.Range("E21").Value = "  +0.41%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "67.66"
$ws.Range("E22").Value = "  +0.03%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "237.09"
$ws.Range("E23").Value = "  -1.21%  "
$ws.Range("E24").Value = "  +2.90%  "
$ws.Range("E25").Value = "  +0.16%  "
$ws.Range("E26").Value = "  +0.33%  "
$ws.Range("E27").Value = "  +0.13%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "166.02"
$ws.Range("E29").Value = "  -0.05%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "9.07"
$ws.Range("E30").Value = "  +0.59%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "32.89"
$ws.Range("E31").Value = "  +0.52%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "18.18"
$ws.Range("E33").Value = "  +7.75%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.99"
$ws.Range("E34").Value = "  -0.05%  "
$ws.Range("E35").Value = "  -7.40%  "
$ws.Range("E36").Value = "  -0.95%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.0691"
$ws.Range("E37").Value = "  +1.21%  "
$ws.Range("E38").Value = "  -0.29%  "
$ws.Range("E39").Value = "  +0.94%  "
$ws.Range("E40").Value = "  +1.65%  "
$ws.Range("E41").Value = "  -0.22%  "
$ws.Range("D42").Value = "2.001.40"
$ws.Range("E42").Value = "  -0.52%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "10.43"
$ws.Range("E43").Value = "  +3.53%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0280"
$ws.Range("E44").Value = "  -0.20%  "
$ws.Range("B45").Value = "EnergySwap"
$ws.Range("C45").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "18.06"
$ws.Range("E45").Value = "  +5.81%  "
$ws.Range("B46").Value = "ApeXProtocol"
$ws.Range("C46").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.12"
$ws.Range("E46").Value = "  +0.01%  "
$ws.Range("E47").Value = "  +0.46%  "
$ws.Range("B48").Value = "MultiversX"
$ws.Range("C48").Value = "https://coinranking.com/coin/omwkOTglq+multiversx-egld"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "53.62"
$ws.Range("E48").Value = "  +1.13%  "
$ws.Range("B49").Value = "RocketPoolETH"
$ws.Range("C49").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D49").Value = "2.532.52"
$ws.Range("E49").Value = "  +0.68%  "
$ws.Range("B50").Value = "HuobiToken"
$ws.Range("C50").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.83"
$ws.Range("E50").Value = "  -3.41%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "71.89"
$ws.Range("E51").Value = "  +0.09%  "
